$wb = $excel.ActiveWorkbook

# "展览" (Exhibition) sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 266
$ws1.Range("F5").Value = 794
$ws1.Range("F6").Value = 1941

# "全部类型" (All types) sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 266
$ws4.Range("F7").Value = 794
$ws4.Range("F8").Value = 1941
